# Salmonella/Coccidia updation
# Replaces the old "TestAutomation3427" cartridge run data (rows 2-13) with
# the new "TestCartridge5813" run: each row gets a (possibly reshuffled)
# Lane number, the new Cartridge ID, a Result ID, and a Lab Sample ID.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lane (col B) holds digit-only text ("1".."12"). A plain .Value assignment
# would be auto-coerced to a number by Excel, so the range is temporarily
# switched to Text format, written, then restored to the workbook's usual
# (unstyled) look so the cells keep rendering/storing as text.
$laneRange = $ws.Range("B2:B13")
$laneRange.NumberFormat = "@"

$rows = @{
    2  = @{ B = "2";  D = "A0789802"; R = "TestAutomation15813U" }
    3  = @{ B = "8";  D = "A0789808"; R = "TestAutomation35813U" }
    4  = @{ B = "12"; D = "A0789812"; R = "TestAutomation45813U" }
    5  = @{ B = "5";  D = "A0789805"; R = "TestAutomation25813U" }
    6  = @{ B = "6";  D = "A0789806"; R = "TestAutomation25813U" }
    7  = @{ B = "10"; D = "A0789810"; R = "TestAutomation45813U" }
    8  = @{ B = "9";  D = "A0789809"; R = "TestAutomation35813U" }
    9  = @{ B = "3";  D = "A0789803"; R = "TestAutomation15813U" }
    10 = @{ B = "11"; D = "A0789811"; R = "TestAutomation45813U" }
    11 = @{ B = "4";  D = "A0789804"; R = "TestAutomation25813U" }
    12 = @{ B = "7";  D = "A0789807"; R = "TestAutomation35813U" }
    13 = @{ B = "1";  D = "A0789801"; R = "TestAutomation15813U" }
}

foreach ($r in 2..13) {
    $vals = $rows[$r]
    $ws.Range("B$r").Value = $vals.B
    $ws.Range("C$r").Value = "TestCartridge5813"
    $ws.Range("D$r").Style = "Normal"
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("R$r").Style = "Normal"
    $ws.Range("R$r").Value = $vals.R
}

# Drop the transient Text number-format now that every Lane cell has been
# written, so the cells end up back on the sheet's normal (unstyled) look.
$laneRange.Style = "Normal"
